# The deck's slide master (theme/theme1.xml, "Integral") is switched back to
# the default "Office Theme" colour scheme (which previously only backed the
# notes master's theme/theme2.xml). Re-create that by pushing the stock
# Office theme RGB values onto the presentation's (slide-master-bound) theme
# colour scheme, one swatch at a time, through the live Slide object so the
# clrScheme's structure (and its 12 named slots) stays intact.

function HexToRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Standard Office default theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
$officeTheme = @("000000","FFFFFF","44546A","E7E6E6","5B9BD5","ED7D31","A5A5A5","FFC000","4472C4","70AD47","0563C1","954F72")

$slide = $p.Slides.Item(1)
$colorScheme = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $swatch = $colorScheme.Colors($i)
    $swatch.RGB = HexToRGB $officeTheme[$i - 1]
}
